$d = $word.ActiveDocument

# The Title, Author and Abstract paragraphs originally have their text
# split across many single-word/single-space runs. Collapse each of
# those paragraphs down into one run holding the full text by running
# Find & Replace (replacing the text with itself) scoped to just that
# paragraph's range -- this merges all the runs inside the range into a
# single run without touching any other part of the document.

$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute(
    "Questions: Vector addition and scalar multiplication",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Questions: Vector addition and scalar multiplication", 2)

$authorRange = $d.Paragraphs(2).Range
$authorRange.Find.Execute(
    "Renee Knapp, Kin Wang Pang",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Renee Knapp, Kin Wang Pang", 2)

$abstractRange = $d.Paragraphs(4).Range
$abstractRange.Find.Execute(
    "A selection of questions for the study guide on vector addition and scalar multiplication.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on vector addition and scalar multiplication.", 2)
